$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "File Name" values in column A with the electrode location
# codes that currently live in column C, for every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 62 }

for ($r = 2; $r -le $lastRow; $r++) {
    $loc = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 1).Value = $loc
}

# Update header text.
$ws.Cells.Item(1, 1).Value = "Loc"
$ws.Cells.Item(1, 2).Value = "P_max"

# Remove the now-redundant "Electrode Locations" column entirely.
$ws.Columns.Item(3).Delete()
